# Refresh the "cryptos" price/volume table with the latest scrape, and
# reorder three rows (43-45) whose rank shuffled: ThetaToken, FirstDigitalUSD,
# Maker -> FirstDigitalUSD, Maker, ThetaToken.
#
# Note: several "Price" values look like plain numbers (e.g. "190.03") but
# must stay TEXT, matching the sheet's existing inline-string cells. A bare
# Range.Value assignment would let Excel auto-coerce those into numeric
# cells, so for any new value that parses as a number we prefix it with a
# leading apostrophe (the same quote-prefix trick a person typing into Excel
# would use) to force it to stay text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.842.20'
$ws.Range('E2').Value = '  +7.86%  '
$ws.Range('D3').Value = '3.512.57'
$ws.Range('E3').Value = '  +10.83%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '''190.03'
$ws.Range('E5').Value = '  +12.83%  '
$ws.Range('D6').Value = '''553.30'
$ws.Range('E6').Value = '  +7.66%  '
$ws.Range('D7').Value = '3.503.31'
$ws.Range('E7').Value = '  +10.58%  '
$ws.Range('D8').Value = '''0.608'
$ws.Range('E8').Value = '  +4.16%  '
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').Value = '''0.636'
$ws.Range('E10').Value = '  +8.04%  '
$ws.Range('D11').Value = '''0.151'
$ws.Range('E11').Value = '  +18.48%  '
$ws.Range('D12').Value = '''55.67'
$ws.Range('E12').Value = '  +7.85%  '
$ws.Range('D13').Value = '''0.0000269'
$ws.Range('E13').Value = '  +8.84%  '
$ws.Range('D14').Value = '''9.42'
$ws.Range('E14').Value = '  +8.11%  '
$ws.Range('D15').Value = '4.066.21'
$ws.Range('E15').Value = '  +9.92%  '
$ws.Range('D16').Value = '3.506.90'
$ws.Range('E16').Value = '  +10.04%  '
$ws.Range('D18').Value = '''18.30'
$ws.Range('E18').Value = '  +8.63%  '
$ws.Range('D19').Value = '66.804.12'
$ws.Range('E19').Value = '  +7.62%  '
$ws.Range('D20').Value = '''11.84'
$ws.Range('E20').Value = '  +10.51%  '
$ws.Range('D21').Value = '''0.996'
$ws.Range('E21').Value = '  +5.96%  '
$ws.Range('D22').Value = '''413.12'
$ws.Range('E22').Value = '  +15.39%  '
$ws.Range('D23').Value = '''3.94'
$ws.Range('E23').Value = '  +7.33%  '
$ws.Range('D24').Value = '''85.41'
$ws.Range('E24').Value = '  +7.73%  '
$ws.Range('E25').Value = '  +11.22%  '
$ws.Range('D26').Value = '''11.14'
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('E27').Value = '  +14.90%  '
$ws.Range('D28').Value = '''6.13'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('E29').Value = '  +9.16%  '
$ws.Range('D30').Value = '''8.89'
$ws.Range('E30').Value = '  +11.23%  '
$ws.Range('D31').Value = '''30.37'
$ws.Range('E31').Value = '  +9.54%  '
$ws.Range('D32').Value = '''655.43'
$ws.Range('E32').Value = '  +2.70%  '
$ws.Range('D33').Value = '''6.73'
$ws.Range('E33').Value = '  +7.48%  '
$ws.Range('D34').Value = '''11.76'
$ws.Range('E34').Value = '  +7.13%  '
$ws.Range('E35').Value = '  +9.30%  '
$ws.Range('D36').Value = '''59.83'
$ws.Range('E36').Value = '  +5.72%  '
$ws.Range('D37').Value = '''38.95'
$ws.Range('E37').Value = '  +9.40%  '
$ws.Range('D38').Value = '0.0₃0811'
$ws.Range('E38').Value = '  +19.35%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('E40').Value = '  +7.90%  '
$ws.Range('E41').Value = '  +13.65%  '
$ws.Range('D42').Value = '''3.35'
$ws.Range('E42').Value = '  +23.49%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''0.998'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '3.002.26'
$ws.Range('E44').Value = '  +5.99%  '
$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').Value = '''2.93'
$ws.Range('E45').Value = '  +17.18%  '
$ws.Range('D46').Value = '''2.64'
$ws.Range('E46').Value = '  +8.66%  '
$ws.Range('D47').Value = '''3.36'
$ws.Range('E47').Value = '  +15.55%  '
$ws.Range('D48').Value = '''0.0418'
$ws.Range('E48').Value = '  +10.73%  '
$ws.Range('E49').Value = '  +5.07%  '
$ws.Range('D50').Value = '''8.97'
$ws.Range('E50').Value = '  +20.73%  '
$ws.Range('E51').Value = '  +8.32%  '
